# Add new power plant types to the Electricity Source subscript on the
# "CRbQ" sheet (issues #280 and #99).
#
# Remember which sheet/cell is selected right now so we can restore it
# (the "About" sheet is the one shown when the workbook is opened) and
# only change view-state on the sheet we are actually editing.
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("CRbQ")

$newPlantTypes = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

# Existing data runs through row 18 (18 rows of plant types under the
# header row) and columns B..AF (years 2020-2050). Append six more rows,
# one per new plant type, with 0 in every year column - matching the
# pattern used by every other plant-type row already on the sheet.
$firstNewRow = 19
$lastCol = 32   # column AF

$ws.Select() | Out-Null

$row = $firstNewRow
foreach ($label in $newPlantTypes) {
    $labelCell = $ws.Cells.Item($row, 1)
    $labelCell.Value = $label
    $labelCell.Font.Bold = $true

    $dataRange = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, $lastCol))
    $dataRange.Value = 0

    $row++
}

# Reflect where the user ended up looking after adding the rows.
$ws.Range("A25").Select() | Out-Null

# Restore the originally active sheet so we don't change which tab is
# marked as selected in the saved workbook.
$about.Select() | Out-Null
